$wb = $excel.ActiveWorkbook

# Apply the "handoff transform failed" update to both locale sheets
# (zh-cn and de-de): the Status moves from "Ready for handoff" to
# "Handoff transform failed", the Latest Handoff File link/value is
# cleared out, the Latest Handoff Datetime and Latest Handback DateTime
# revert to the zero-date sentinel, and the Handoff Reason flips from
# "Include" to "Ignored".
$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the hyperlink attached to the "Latest Handoff File" cell (C2)
    # before clearing the cell itself.
    $linksToRemove = @()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$2') {
            $linksToRemove += $hl
        }
    }
    foreach ($hl in $linksToRemove) {
        $hl.Delete()
    }

    $ws.Range("B2").Value = "Handoff transform failed"
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Ignored"
}
